$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the extra "Peter" rows (rows 5-12), keeping row 4 to be overwritten
$ws.Range("A5:C12").EntireRow.Delete() | Out-Null

# Row 2 becomes Bala / bala123 / Invalid credentials
$ws.Range("A2").Value = "Bala"
$ws.Range("B2").Value = "bala123"
$ws.Range("C2").Value = "Invalid credentials"

# Row 4 becomes Mark / Mark123 / Invalid credentials
$ws.Range("A4").Value = "Mark"
$ws.Range("B4").Value = "Mark123"
$ws.Range("C4").Value = "Invalid credentials"

# Update the selection to match the target state
$ws.Range("B2").Select() | Out-Null
